$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.83 = 6709.38 pesos`n✅ 6709.38 pesos = 1.83 = 941.31 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 546.995
$ws2.Range("O10").Value = 3670
$ws2.Range("N12").Value = 3662.76
$ws2.Range("O12").Value = 513.877
